$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "PASS"

$ws.Range("C6").Value = "xfsimmvonglnkvk@gmail.com"
$ws.Range("D6").Value = "opqfmLGASI5"
$ws.Range("E6").Value = "pass"
